$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("C14").Value = 250.56
$wsVentasGrupo.Range("L14").Value = 3217
$wsVentasGrupo.Range("M14").Value = 1081.17
$wsVentasGrupo.Range("C29").Value = "1 de 27"
$wsVentasGrupo.Range("M29").Value = "3 de 27"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F14").Value = 3990.41
$wsVentaMensual.Range("F29").Value = 14727.56

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column E gets one unit wider (raw OOXML width 22 -> 23).
# Excel's ColumnWidth property stores a value 5/6 below the raw <col width>
# attribute for this workbook's default font, so back that offset out here.
$wsCumplimiento.Columns.Item(5).ColumnWidth = 22.166666666666668

$wsCumplimiento.Range("D2").Value = 250.56
$wsCumplimiento.Range("E2").Value = 93.72460462948601
$wsCumplimiento.Range("F2").Value = 0.7277699805068222

$wsCumplimiento.Range("D15").Value = 3526.47
$wsCumplimiento.Range("E15").Value = -1888.47
$wsCumplimiento.Range("F15").Value = 2.152912087912088

$wsCumplimiento.Range("D16").Value = 4994.29
$wsCumplimiento.Range("E16").Value = 12091.6
$wsCumplimiento.Range("F16").Value = 0.2923049369977215

$wsCumplimiento.Range("D19").Value = 14962.03
$wsCumplimiento.Range("E19").Value = 12219.28093005039
$wsCumplimiento.Range("F19").Value = 0.5504528474915711
